$wb = $excel.ActiveWorkbook

# Sheet1 is "card" (sheet1.xml), Sheet2 is "Sheet1" (sheet2.xml)
$wsCard = $wb.Worksheets.Item(1)
$wsSheet1 = $wb.Worksheets.Item(2)

# Add the two new rows of data to the "Sheet1" worksheet (sheet2.xml): row 9
# A9 = "*Browser", B9 = "IE11/FF/Chrome" with the same green "passed" fill style as B2:B7
$wsSheet1.Range("A9").Value = "*Browser"
$wsSheet1.Range("B9").Value = "IE11/FF/Chrome"
$wsSheet1.Range("B9").Interior.Color = $wsSheet1.Range("B2").Interior.Color

# Update the selection on Sheet1 (sheet2.xml) to the full used range A1:F9
$wsSheet1.Range("A1:F9").Select()

# Make "Sheet1" (sheet2.xml) the active/selected tab, as in the diff
# (tabSelected moves from "card" to "Sheet1", and workbookView activeTab becomes 1)
$wsSheet1.Activate()
